# Update Name of Algo
# Applies the data corrections to result_data_RandomForest.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 12.0688
$ws.Range("B8").Value = 4.816899999999996
$ws.Range("B10").Value = 8.607100000000003
$ws.Range("B12").Value = 6.075799999999999
$ws.Range("C13").Value = -12.24979999999999
$ws.Range("B18").Value = 5.342100000000003
$ws.Range("E20").Value = 13.32249999999999
$ws.Range("B25").Value = 5.409999999999993
